# identaçao, comentários, novas MF simétricas
# Update CLPVariation_pred (N) and erro_CLP (O) columns with values recomputed
# from the new symmetric membership functions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => (N value, O value)
$updates = @{
    2  = @(0.6562615384615387, 0.07373846153846131)
    3  = @(-0.8105953846153846, 0.009404615384615345)
    4  = @(0.02046153846153852, 0.09953846153846148)
    6  = @(0.5500416666666665, 0.05004166666666654)
    9  = @(-0.5399939130434782, 0.1100060869565218)
    11 = @(0.9070000000000001, 0.05700000000000016)
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Range("N$row").Value = $values[0]
    $ws.Range("O$row").Value = $values[1]
}

$wb.Save()
